# Day2_Communication.pptx restructuring:
#  - insert a new icebreaker slide before the existing title slide
#  - append 9 more new slides (section dividers + "5 Minute Break" slides +
#    a "Group Project" slide) after the existing title slide
#  - group the resulting 11 slides into 7 named sections

$p = $ppt.ActivePresentation

# EMU -> Point conversion helper (PowerPoint's object model positions shapes
# in points; 1 point = 12700 EMU).
$EMU = 12700

# ---------------------------------------------------------------------
# 1. New slide 1: icebreaker question (Title Only layout), inserted before
#    the existing (only) slide.
# ---------------------------------------------------------------------
$s1 = $p.Slides.Add(1, 3)
$t1 = $s1.Shapes.Item(1)
$t1.TextFrame.TextRange.Text = "Without modern medical intervention, how many times would you have died so far?"
$t1.Left = 838200 / $EMU
$t1.Top = 2766218 / $EMU
$t1.Width = 10515600 / $EMU
$t1.Height = 1325563 / $EMU

# ---------------------------------------------------------------------
# The pre-existing slide (title slide, "Communicating Research with
# Effective Storytelling") is now at position 2 and needs no changes.
# ---------------------------------------------------------------------

# ---------------------------------------------------------------------
# 2. Slide 3: empty "Storytelling" section divider (Title Only layout).
# ---------------------------------------------------------------------
$s3 = $p.Slides.Add(3, 3)

# ---------------------------------------------------------------------
# 3. Slide 4: "5 Minute Break" (Title Only layout, centered).
# ---------------------------------------------------------------------
$s4 = $p.Slides.Add(4, 3)
$t4 = $s4.Shapes.Item(1)
$t4.TextFrame.TextRange.Text = "5 Minute Break"
$t4.TextFrame.TextRange.ParagraphFormat.Alignment = 2
$t4.Left = 838200 / $EMU
$t4.Top = 2766219 / $EMU
$t4.Width = 10515600 / $EMU
$t4.Height = 1325563 / $EMU

# ---------------------------------------------------------------------
# 4. Slide 5: empty "PowerPoint Aesthetics" section divider (Title and
#    Content layout).
# ---------------------------------------------------------------------
$s5 = $p.Slides.Add(5, 2)

# ---------------------------------------------------------------------
# 5. Slide 6: "5 Minute Break" (Title Only layout, centered).
# ---------------------------------------------------------------------
$s6 = $p.Slides.Add(6, 3)
$t6 = $s6.Shapes.Item(1)
$t6.TextFrame.TextRange.Text = "5 Minute Break"
$t6.TextFrame.TextRange.ParagraphFormat.Alignment = 2
$t6.Left = 838200 / $EMU
$t6.Top = 2766219 / $EMU
$t6.Width = 10515600 / $EMU
$t6.Height = 1325563 / $EMU

# ---------------------------------------------------------------------
# 6. Slide 7: empty "Morph" section divider (Title and Content layout).
# ---------------------------------------------------------------------
$s7 = $p.Slides.Add(7, 2)

# ---------------------------------------------------------------------
# 7. Slide 8: "5 Minute Break" (Title Only layout, centered).
# ---------------------------------------------------------------------
$s8 = $p.Slides.Add(8, 3)
$t8 = $s8.Shapes.Item(1)
$t8.TextFrame.TextRange.Text = "5 Minute Break"
$t8.TextFrame.TextRange.ParagraphFormat.Alignment = 2
$t8.Left = 838200 / $EMU
$t8.Top = 2766219 / $EMU
$t8.Width = 10515600 / $EMU
$t8.Height = 1325563 / $EMU

# ---------------------------------------------------------------------
# 8. Slide 9: empty "Tips & Tricks" section divider (Title and Content
#    layout).
# ---------------------------------------------------------------------
$s9 = $p.Slides.Add(9, 2)

# ---------------------------------------------------------------------
# 9. Slide 10: "5 Minute Break" (Title Only layout, centered).
# ---------------------------------------------------------------------
$s10 = $p.Slides.Add(10, 3)
$t10 = $s10.Shapes.Item(1)
$t10.TextFrame.TextRange.Text = "5 Minute Break"
$t10.TextFrame.TextRange.ParagraphFormat.Alignment = 2
$t10.Left = 838200 / $EMU
$t10.Top = 2766219 / $EMU
$t10.Width = 10515600 / $EMU
$t10.Height = 1325563 / $EMU

# ---------------------------------------------------------------------
# 10. Slide 11: "Group Project" (Title and Content layout).
# ---------------------------------------------------------------------
$s11 = $p.Slides.Add(11, 2)
$t11 = $s11.Shapes.Item(1)
$t11.TextFrame.TextRange.Text = "Group Project"

# ---------------------------------------------------------------------
# 11. Group the 11 slides into the 7 named sections.
# ---------------------------------------------------------------------
$sp = $p.SectionProperties
$sp.AddBeforeSlide(1, "Personal intros + icebreaker")
$sp.AddBeforeSlide(2, "Intro")
$sp.AddBeforeSlide(3, "Storytelling")
$sp.AddBeforeSlide(5, "PowerPoint Aesthetics")
$sp.AddBeforeSlide(7, "Morph")
$sp.AddBeforeSlide(9, "Tips & Tricks")
$sp.AddBeforeSlide(11, "Project")
